$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StringLocalizations_BasicText")

# Insert two new rows at position 34 (pushes the "Turns Survived" block and
# everything after it down by two rows) and populate them with the new
# "Cases Closed Positively/Negatively" localization strings.
$ws.Rows.Item(34).Insert()
$ws.Rows.Item(35).Insert()

$ws.Cells.Item(34, 1).Value = "BASIC_TEXT_CASES_CLOSED_POSITIVE"
$ws.Cells.Item(34, 2).Value = "Cases Closed Positively"
$ws.Cells.Item(34, 3).Value = "Cases Closed Positively"
$ws.Cells.Item(34, 4).Value = "XXXX"
$ws.Cells.Item(34, 5).Value = "Casos cerrados positivamente"

$ws.Cells.Item(35, 1).Value = "BASIC_TEXT_CASES_CLOSED_NEGATIVE"
$ws.Cells.Item(35, 2).Value = "Cases Closed Negatively"
$ws.Cells.Item(35, 3).Value = "Cases Closed Negatively"
$ws.Cells.Item(35, 4).Value = "XXXX"
$ws.Cells.Item(35, 5).Value = "Casos cerrados negativamente"

# Reflect the author's new scroll position / selection on the sheet.
[void]$ws.Activate()
$ws.Range("C25").Select() | Out-Null
